# Adds a new numbered step after "Mendownload file gith di internet."
# describing what to do once the git installer has finished downloading.

$d = $word.ActiveDocument

# Locate the existing numbered paragraph that ends the download step. We
# search instead of hard-coding a paragraph index so the script still works
# if the surrounding content shifts slightly.
$anchorText = "Mendownload file gith di internet."
$searchRange = $d.Content
$found = $searchRange.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph '$anchorText'"
}

$anchorPara = $searchRange.Paragraphs(1)

# Insert a brand-new paragraph right after it. Word automatically carries
# over the anchor paragraph's formatting (Heading1 numbered-list style,
# numId 1 / ilvl 0, spacing, indent and run font/bold/color), matching the
# existing list item's look.
$anchorPara.Range.InsertParagraphAfter()
$newPara = $anchorPara.Next()

$newPara.Range.Text = "Setelah terdownload, install file git tadi."
